$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "29.381.00"
$ws.Range("E2").Value = "  +0.46%  "
Set-TextValue $ws.Range("D3") "1.874.23"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.03%  "
Set-TextValue $ws.Range("D5") "0.7114"
$ws.Range("E5").Value = "  -0.50%  "
Set-TextValue $ws.Range("D6") "242.01"
$ws.Range("E6").Value = "  +0.49%  "
$ws.Range("E7").Value = "  -0.01%  "
Set-TextValue $ws.Range("D8") "0.07792"
$ws.Range("E8").Value = "  +0.87%  "
$ws.Range("E9").Value = "  +0.73%  "
$ws.Range("E10").Value = "  +0.80%  "
Set-TextValue $ws.Range("D11") "0.08451"
$ws.Range("E11").Value = "  +1.60%  "
Set-TextValue $ws.Range("D12") "1.870.79"
$ws.Range("E12").Value = "  +0.27%  "
Set-TextValue $ws.Range("D13") "5.236"
$ws.Range("E13").Value = "  +0.30%  "
Set-TextValue $ws.Range("D14") "0.7117"
$ws.Range("E14").Value = "  -0.83%  "
Set-TextValue $ws.Range("D15") "91.14"
$ws.Range("E15").Value = "  +0.28%  "
Set-TextValue $ws.Range("D16") "29.384.06"
$ws.Range("E16").Value = "  +0.47%  "
Set-TextValue $ws.Range("D17") "6.047"
$ws.Range("E17").Value = "  +0.85%  "
Set-TextValue $ws.Range("D18") "0.000008221"
$ws.Range("E18").Value = "  +5.21%  "
Set-TextValue $ws.Range("D19") "241.05"
$ws.Range("E19").Value = "  -0.99%  "
Set-TextValue $ws.Range("D20") "13.25"
$ws.Range("E20").Value = "  +0.59%  "
Set-TextValue $ws.Range("D21") "2.120.29"
$ws.Range("E21").Value = "  -0.03%  "
Set-TextValue $ws.Range("D22") "0.9999"
Set-TextValue $ws.Range("D23") "7.778"
$ws.Range("E23").Value = "  -1.83%  "
$ws.Range("E24").Value = "  -0.01%  "
Set-TextValue $ws.Range("D25") "0.1606"
$ws.Range("E25").Value = "  -0.46%  "
Set-TextValue $ws.Range("D26") "163.55"
$ws.Range("E26").Value = "  +0.47%  "
Set-TextValue $ws.Range("D27") "9.045"
$ws.Range("E27").Value = "  +1.48%  "
Set-TextValue $ws.Range("D28") "18.49"
$ws.Range("E28").Value = "  -0.61%  "
Set-TextValue $ws.Range("D29") "1.511"
$ws.Range("E29").Value = "  +0.78%  "
Set-TextValue $ws.Range("D30") "4.432"
$ws.Range("E30").Value = "  -0.01%  "
Set-TextValue $ws.Range("D31") "1.290"
$ws.Range("E31").Value = "  -5.12%  "
Set-TextValue $ws.Range("D32") "4.300"
$ws.Range("E32").Value = "  +0.82%  "
Set-TextValue $ws.Range("D33") "0.05276"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +0.31%  "
Set-TextValue $ws.Range("D36") "0.7451"
$ws.Range("E36").Value = "  -8.59%  "
Set-TextValue $ws.Range("D37") "2.696"
$ws.Range("E37").Value = "  +0.52%  "
$ws.Range("E38").Value = "  +0.49%  "
Set-TextValue $ws.Range("D39") "1.214.40"
$ws.Range("E39").Value = "  +4.82%  "
Set-TextValue $ws.Range("D40") "2.721"
$ws.Range("E40").Value = "  +1.03%  "
Set-TextValue $ws.Range("D41") "6.483"
$ws.Range("E41").Value = "  +4.35%  "
Set-TextValue $ws.Range("D42") "0.8889"
$ws.Range("E42").Value = "  -0.68%  "
Set-TextValue $ws.Range("D43") "72.72"
$ws.Range("E43").Value = "  -0.11%  "
$ws.Range("E44").Value = "  +6.96%  "
$ws.Range("E45").Value = "  +0.01%  "
Set-TextValue $ws.Range("D46") "2.017.95"
$ws.Range("E46").Value = "  -0.86%  "
Set-TextValue $ws.Range("D47") "1.814"
$ws.Range("E47").Value = "  +1.72%  "
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws.Range("D48") "0.5210"
$ws.Range("E48").Value = "  +0.63%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D49") "9.365"
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue $ws.Range("D50") "0.4319"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("B51").Value = "Aptos"
$ws.Range("C51").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D51") "7.088"
$ws.Range("E51").Value = "  +0.27%  "
